# Add new command "dragAndDrop(fromLocator,toLocator)" to the hidden
# '#system' sheet's "web" command list (column U), keeping the list in
# alphabetical order, then grow the "web" named range to cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# The new entry belongs right before "editLocalStorage(key,value)", which
# currently lives at U59 -- i.e. between doubleClickByLabelAndWait(...) at
# U58 and editLocalStorage(...) at U59.
$insertRow = 59
$lastRow = 111
$newLastRow = $lastRow + 1

# Capture the existing values for U59:U111 (column U only -- column E, which
# holds an unrelated "desktop" command list, must stay untouched).
$vals = @()
for ($r = $insertRow; $r -le $lastRow; $r++) {
    $vals += $ws.Range("U$r").Value2
}

# Shift them all down by one row, from the bottom up isn't required here
# since we already captured everything above before overwriting.
for ($i = 0; $i -lt $vals.Length; $i++) {
    $targetRow = $insertRow + 1 + $i
    $ws.Range("U$targetRow").Value = $vals[$i]
}

# Place the new command in the freed-up slot.
$ws.Range("U$insertRow").Value = "dragAndDrop(fromLocator,toLocator)"

# Grow the "web" defined name so it covers the newly added row.
$wb.Names.Item("web").RefersTo = "='#system'!`$U`$2:`$U`$$newLastRow"
